# OpenDataServices/flatten-tool#100 - fix org location title headings in template
# The "fun_location" (funding organization location) and "rec_location"
# (recipient organization location) summary-table sheets both used the
# generic "Location:..." column headings. Prefix them with the owning
# organization so it's clear which org the location block belongs to.

$wb = $excel.ActiveWorkbook

$locationColumns = @("Identifier", "Name", "Country Code", "Latitude", "Longitude", "Description", "Geographic Code", "Geographic Code Type", "Last modified")
$columnLetters = @("C", "D", "E", "F", "G", "H", "I", "J", "K")

$fundingLocationSheet = $wb.Worksheets.Item("fun_location")
for ($i = 0; $i -lt $locationColumns.Length; $i++) {
    $cellRef = $columnLetters[$i] + "1"
    $fundingLocationSheet.Range($cellRef).Value = "Funding Org:Location:" + $locationColumns[$i]
}

$recipientLocationSheet = $wb.Worksheets.Item("rec_location")
for ($i = 0; $i -lt $locationColumns.Length; $i++) {
    $cellRef = $columnLetters[$i] + "1"
    $recipientLocationSheet.Range($cellRef).Value = "Recipient Org:Location:" + $locationColumns[$i]
}
